$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item(2)
$new = $wb.Worksheets.Add($null, $ws2)
$new.Range("B3").NumberFormat = "h:mm"
$new.Range("B3").Value = 0.52083333333333337
$new.Range("B4").NumberFormat = "h:mm"
$new.Range("B4").Value = 0.5625
